$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Handoff transform failed" -> "Ready for handoff" (shared string reused across sheets;
# updating every cell that showed the old text keeps all of them in sync)
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("B2").Value = "Ready for handoff"

# zh-cn row 2: new "Latest Handoff File" hyperlink + datetime + reason
$zhFile = "a7121821-1a0d-4dfc-9a39-c640b6860ad5.6643f4fe47d4d0b1c828ab9b15c1633a14e49f6d.zh-cn.xlf"
$zhUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8daa92ec74ac0976e3165cf36a85b0787a765527/e2e/" + $zhFile
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhUrl, "", "", $zhFile)
$wsZhCn.Range("D2").Value = "2016-01-25 14:02:43"
$wsZhCn.Range("H2").Value = "Include"

# de-de row 2: new "Latest Handoff File" hyperlink + datetime + reason
$deFile = "a7121821-1a0d-4dfc-9a39-c640b6860ad5.6643f4fe47d4d0b1c828ab9b15c1633a14e49f6d.de-de.xlf"
$deUrl = "https://github.com/OpenLocalizationTest/oltest/blob/8daa92ec74ac0976e3165cf36a85b0787a765527/e2e/" + $deFile
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deUrl, "", "", $deFile)
$wsDeDe.Range("D2").Value = "2016-01-25 14:02:52"
$wsDeDe.Range("H2").Value = "Include"
